$p = $ppt.ActivePresentation

# --- 1. Update the "datetimeFigureOut" date placeholder text on every slide
#        layout of the (single) slide master: "9/28/2018" -> "12/7/2018"
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.CustomLayouts.Count; $i++) {
    $cl = $sm.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $shp = $cl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/28/2018") {
                $tr.Text = "12/7/2018"
            }
        }
    }
}

# --- 2. Update the same date field on the Notes Master: "28/9/2018" -> "7/12/2018"
$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $shp = $nm.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "28/9/2018") {
            $tr.Text = "7/12/2018"
        }
    }
}

# --- 3. Rename the "Karma" testing-tool label to "Jest" on slide 1
$s = $p.Slides.Item(1)
for ($k = 1; $k -le $s.Shapes.Count; $k++) {
    $shp = $s.Shapes.Item($k)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Karma") {
            $shp.TextFrame.TextRange.Text = "Jest"
        }
    }
}
